$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking export stores Price/Volume as plain text. Force text format on the
# Price cells we touch so numeric-looking strings like "597.70" or "14.50" keep
# their trailing zero instead of being coerced into a Number by COM.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.860.69"
$ws.Range("E2").Value = "  -1.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.132.58"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.70"
$ws.Range("E5").Value = "  -2.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.01"
$ws.Range("E6").Value = "  -4.50%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.127.18"
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.34"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000252"
$ws.Range("E13").Value = "  -2.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.45"
$ws.Range("E14").Value = "  -3.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.663.74"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("E16").Value = "  +2.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.839.66"
$ws.Range("E17").Value = "  -1.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.139.60"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.75"
$ws.Range("E19").Value = "  -2.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "482.13"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.50"
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.704"
$ws.Range("E22").Value = "  -2.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.67"
$ws.Range("E23").Value = "  -3.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.72"
$ws.Range("E24").Value = "  +4.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.11"
$ws.Range("E25").Value = "  -4.69%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.74"
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.14"
$ws.Range("E28").Value = "  -7.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.95"
$ws.Range("E29").Value = "  -3.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.04"
$ws.Range("E30").Value = "  -3.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.24"
$ws.Range("E31").Value = "  +2.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.111"
$ws.Range("E32").Value = "  -7.46%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.59"
$ws.Range("E34").Value = "  -4.09%  "
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.01"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.63"
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0726"
$ws.Range("E38").Value = "  -8.56%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0395"
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.87"
$ws.Range("E40").Value = "  -10.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "427.43"
$ws.Range("E41").Value = "  -7.66%  "
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.29"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.904.76"
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.258"
$ws.Range("E45").Value = "  -4.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.16"
$ws.Range("E46").Value = "  -7.46%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.36"
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.114"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.54"
$ws.Range("E50").Value = "  -4.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.45"
$ws.Range("E51").Value = "  +0.49%  "
